$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (AIRCRAFT F1 total) ---
$ws.Range("C2").Value = 0.44290000000000002
$ws.Range("D2").Value = 0.25840000000000002
$ws.Range("E2").Value = 0.32519999999999999

# --- Row 3 (1.0 LEGACY ABS F1) ---
$ws.Range("C3").Value = 1.9422999999999999
$ws.Range("D3").Value = 1.8139000000000001
# E3 unchanged (2.8537)

# --- Row 4 (1L EETC F1) ---
$ws.Range("C4").Value = -0.021600000000000001
$ws.Range("D4").Value = -0.021600000000000001
$ws.Range("E4").Value = -0.032899999999999999

# --- Row 5 (2L EETC F1) ---
$ws.Range("C5").Value = 0.0177
# D5 unchanged (0.0176)
$ws.Range("E5").ClearContents()

# --- Row 6 (3.0 MEZZ ABS F1) ---
$ws.Range("C6").Value = 0.58979999999999999
$ws.Range("D6").Value = 0.37759999999999999
$ws.Range("E6").Value = 0.55469999999999997

# --- Row 7 (3.0 SENIOR ABS F1) ---
$ws.Range("C7").Value = 0.1978
$ws.Range("D7").Value = 0.1164
$ws.Range("E7").Value = 0.10920000000000001

# --- Row 8 (AIR UNSECURED F1) ---
$ws.Range("C8").Value = 0.2298
$ws.Range("D8").Value = 0.13370000000000001
$ws.Range("E8").Value = 0.070999999999999994

# --- Row 9 (AIRCRAFT F1_INCOME) ---
# C9 unchanged (0.0456), D9 unchanged (0.0456)
$ws.Range("E9").ClearContents()

# --- Row 10 (TRADABLE E NOTES F1) ---
$ws.Range("C10").Value = 0.22109999999999999
$ws.Range("D10").Value = 0.12959999999999999
$ws.Range("E10").Value = 0.20030000000000001

# --- Row 11 (CMBS F1 total) ---
$ws.Range("C11").Value = 0.20799999999999999
$ws.Range("D11").Value = 0.10920000000000001
$ws.Range("E11").Value = 0.1457

# --- Row 12 (CMBS 2.0/3.0 IG F1) ---
$ws.Range("C12").Value = 0.2029
$ws.Range("D12").Value = 0.093100000000000002
$ws.Range("E12").Value = 0.1172

# --- Row 13 (CMBS 2.0/3.0 NON-IG F1) ---
$ws.Range("C13").Value = 0.27510000000000001
$ws.Range("D13").Value = 0.15909999999999999
$ws.Range("E13").Value = 0.316

# --- Row 14 (CMBS AGENCY F1) ---
$ws.Range("C14").Value = 0.045999999999999999
$ws.Range("D14").Value = 0.045999999999999999
$ws.Range("E14").Value = 0.041000000000000002

# --- Row 15 (CMBS IO F1) ---
$ws.Range("C15").Value = 0.24199999999999999
$ws.Range("D15").Value = 0.24199999999999999
$ws.Range("E15").Value = 0.45729999999999998

# --- Row 16 (CMBS PRIVATE LOANS) ---
$ws.Range("C16").Value = 0.075200000000000003
$ws.Range("D16").Value = 0.075200000000000003
$ws.Range("E16").Value = 0.1225

# --- Row 17 (CMBS SASB F1) ---
$ws.Range("C17").Value = 0.20250000000000001
$ws.Range("D17").Value = 0.13769999999999999
$ws.Range("E17").Value = 0.1479

# --- Row 18 (CMBS SASB F1_INCOME) ---
$ws.Range("C18").Value = 0.065899999999999997
$ws.Range("D18").Value = 0.065899999999999997
$ws.Range("E18").Value = 0.051900000000000002

# Rows 19-20 (SHORT TERM) are unchanged.

# --- Row 21 (CLO F1 total) ---
$ws.Range("C21").Value = 0.0085000000000000006
$ws.Range("D21").Value = 0.0085000000000000006
$ws.Range("E21").Value = 0.0085000000000000006

# --- Row 22 (CLO AAA ETF F1) ---
$ws.Range("C22").Value = 0.0061999999999999998
$ws.Range("D22").Value = 0.0061999999999999998
$ws.Range("E22").Value = 0.0061999999999999998

# --- Insert two new CLO F1 substrategy rows after row 22: CLO MEZZ F1 / TRUPS MEZZ F1 ---
$ws.Rows("23:24").Insert() | Out-Null

$ws.Range("A23").Value = "CLO F1"
$ws.Range("B23").Value = "CLO MEZZ F1"
$ws.Range("C23").Value = 0.0089999999999999993
$ws.Range("D23").Value = 0.0089999999999999993
$ws.Range("E23").Value = 0.0089999999999999993

$ws.Range("A24").Value = "CLO F1"
$ws.Range("B24").Value = "TRUPS MEZZ F1"
$ws.Range("C24").Value = 0.0089999999999999993
$ws.Range("D24").Value = 0.0089999999999999993
$ws.Range("E24").Value = 0.0089999999999999993

# --- Row 25 (ABS F1 total, formerly row 23) ---
# C25 unchanged (0.077)
$ws.Range("D25").Value = 0.029700000000000001
$ws.Range("E25").Value = 0.085500000000000007

# Rows 26-28 (MEZZ HOME IMPROVEMENT F1 / MEZZ MPL / SENIOR MPL, formerly 24-26) are unchanged.

# --- Update selection to match latest cursor position ---
$ws.Range("F20").Select() | Out-Null
